$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Remove the old header row; this shifts every row up by one (row2->row1, ... row18->row17)
$ws.Rows(1).Delete()

# --- Column A: rename "Method" header + re-letter/re-word the method descriptions ---
$ws.Range("A1").Value = "Extreme Low Flow Method"
$ws.Range("A2").Value = "A. From tree rings back to 1400 AD"
$ws.Range("A3").Value = "B. Collaborator choices in immersive modeling sessions"
$ws.Range("A4").Value = "C.Low Lake Powell releases + gains through Grand Canyon"
$ws.Range("A5").Value = "D. 85%, 65%, and 50% of 2000 to 2018 flow"
$ws.Range("A6").Value = "E. Reclamation's Post 2026 web tool"
$ws.Range("A7").Value = "F. Extreme low flows from Reclamation's ensembles and traces"

# --- Column B: tweak a few Flow Location values ---
$ws.Range("B3").Value = "Lee Ferry"
$ws.Range("B4").Value = "Lake Mead"
$ws.Range("B5").Value = "Lake Powell"
$ws.Range("B7").Value = "Lee Ferry"

# --- Column C: row 7 (previously blank) now has a Flow Type ---
$ws.Range("C7").Value = "Natural"

# --- New column F: "Strategy to Stabilize Lake Levels" ---
$ws.Range("F1").Value = "Strategy to Stabilize Lake Levels"
$ws.Range("F2").Value = "Cap depletions to 10-year lookback period of flow."
$ws.Range("F3").Value = "Divide inflow; Users consume and conserve within their account balance."
$ws.Range("F4").Value = "Rule curve; Consumption equals or less than inflow minus evaporation."
$ws.Range("F5").Value = "Release 95% of regulated inflow."
$ws.Range("F6").Value = "Release to prevent drawdown to 3,490 feet."
$ws.Range("F8").Value = "Immersive modeling in progress"
$ws.Range("F9").Value = "Immersive modeling in progress"
$ws.Range("F10").Value = "Immersive modeling in progress"

# Copy formatting for the new column F cells from the equivalent cells in column A
# (header style for F1, body style for the data rows) so styles are reused rather
# than duplicated.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("F2:F6").PasteSpecial(-4122)
$ws.Range("F8:F10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Column F width
$ws.Columns("F").ColumnWidth = 23.92

# Restore the active selection to match the final state of the edit
$null = $ws.Range("B6").Select()
